# Update crypto price/volume data per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force text interpretation so numeric-looking strings (e.g. "601.51")
    # are not silently coerced into numbers, then restore the default style
    # so no residual formatting difference is left on the cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextCell 'D2' '66.812.91'
Set-TextCell 'E2' '  +3.09%  '
Set-TextCell 'D3' '3.207.61'
Set-TextCell 'E3' '  +2.01%  '
Set-TextCell 'E4' '  +0.06%  '
Set-TextCell 'D5' '601.51'
Set-TextCell 'E5' '  +3.51%  '
Set-TextCell 'D6' '158.13'
Set-TextCell 'E6' '  +7.19%  '
Set-TextCell 'E7' '  -0.03%  '
Set-TextCell 'D8' '3.210.31'
Set-TextCell 'E8' '  +2.09%  '
Set-TextCell 'D9' '0.553'
Set-TextCell 'E9' '  +5.21%  '
Set-TextCell 'D10' '0.160'
Set-TextCell 'E10' '  +1.84%  '
Set-TextCell 'D11' '6.01'
Set-TextCell 'E11' '  -2.49%  '
Set-TextCell 'D12' '0.515'
Set-TextCell 'E12' '  +3.43%  '
Set-TextCell 'E13' '  +1.51%  '
Set-TextCell 'D14' '39.27'
Set-TextCell 'E14' '  +5.72%  '
Set-TextCell 'D15' '3.737.33'
Set-TextCell 'E15' '  +2.31%  '
Set-TextCell 'D16' '66.822.69'
Set-TextCell 'E16' '  +3.08%  '
Set-TextCell 'D17' '7.46'
Set-TextCell 'E17' '  +4.48%  '
Set-TextCell 'D18' '3.211.00'
Set-TextCell 'E18' '  +2.32%  '
Set-TextCell 'E19' '  +0.90%  '
Set-TextCell 'D20' '518.64'
Set-TextCell 'E20' '  +3.78%  '
Set-TextCell 'D21' '15.42'
Set-TextCell 'E21' '  +0.47%  '
Set-TextCell 'D22' '0.742'
Set-TextCell 'E22' '  +4.08%  '
Set-TextCell 'D23' '8.19'
Set-TextCell 'E23' '  +5.43%  '
Set-TextCell 'D24' '15.05'
Set-TextCell 'E24' '  +0.31%  '
Set-TextCell 'D25' '85.33'
Set-TextCell 'E25' '  +1.30%  '
Set-TextCell 'E26' '  -0.05%  '
Set-TextCell 'E27' '  +2.57%  '
Set-TextCell 'E28' '  +3.89%  '
Set-TextCell 'D29' '2.43'
Set-TextCell 'E29' '  +11.01%  '
Set-TextCell 'E30' '  +10.30%  '
Set-TextCell 'D31' '7.07'
Set-TextCell 'E31' '  +10.59%  '
Set-TextCell 'D32' '28.31'
Set-TextCell 'E32' '  +2.75%  '
Set-TextCell 'E33' '  +1.82%  '
Set-TextCell 'E34' '  +0.15%  '
Set-TextCell 'E35' '  +2.35%  '
Set-TextCell 'D36' '527.45'
Set-TextCell 'E36' '  +12.19%  '
Set-TextCell 'D37' '55.05'
Set-TextCell 'E37' '  +0.35%  '
Set-TextCell 'D38' '0.0903'
Set-TextCell 'E38' '  +1.05%  '
Set-TextCell 'D39' '0.0425'
Set-TextCell 'E39' '  +2.05%  '
Set-TextCell 'E40' '  +9.46%  '
Set-TextCell 'D41' '2.95'
Set-TextCell 'E41' '  +1.18%  '
Set-TextCell 'D42' '8.93'
Set-TextCell 'E42' '  +2.23%  '
Set-TextCell 'B43' 'TheGraph'
Set-TextCell 'C43' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell 'D43' '0.308'
Set-TextCell 'E43' '  +8.95%  '
Set-TextCell 'B44' 'PEPE'
Set-TextCell 'C44' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell 'D44' '0.0₃0693'
Set-TextCell 'E44' '  +15.39%  '
Set-TextCell 'E45' '  +3.06%  '
Set-TextCell 'D46' '2.892.94'
Set-TextCell 'E46' '  -2.87%  '
Set-TextCell 'D47' '28.85'
Set-TextCell 'E47' '  +2.11%  '
Set-TextCell 'E48' '  +8.00%  '
Set-TextCell 'E49' '  +3.04%  '
Set-TextCell 'D51' '2.69'
Set-TextCell 'E51' '  +10.57%  '
